# Edit 1: rewrite the "ОСНОВНЫЕ ФИНАНСОВЫЕ ПОКАЗАТЕЛИ КЛИЕНТА (...)" heading paragraph
# so the sample years in the italic hint text go from 2014/2015 to 2016/2017,
# reproducing the exact run layout produced by Word's incremental retyping
# (prefix run, "201" run, single "6" run, "г./3 кв. 201" run, single "7" run, "г.)" run).
$d = $word.ActiveDocument

$headingRange = $d.Content
$foundHeading = $headingRange.Find.Execute("ОСНОВНЫЕ ФИНАНСОВЫЕ ПОКАЗАТЕЛИ КЛИЕНТА", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundHeading) {
    $headingPara = $headingRange.Paragraphs(1)
    $headingParaRange = $headingPara.Range
    $newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:before="0" w:after="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/></w:rPr><w:t>ОСНОВНЫЕ ФИНАНСОВЫЕ ПОКАЗАТЕЛИ КЛИЕНТА (</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i/></w:rPr><w:t>указываются периоды предоставленной отчетности, например, 201</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i/></w:rPr><w:t>6</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i/></w:rPr><w:t>г./3 кв. 201</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i/></w:rPr><w:t>7</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:i/></w:rPr><w:t>г.)</w:t></w:r></w:p>'
    $headingParaRange.InsertXML($newParaXml)

    # Re-fetch the (now replaced) paragraph and make sure the "space before"
    # stays explicit at 0 (InsertXML can fold a default value away).
    $headingRange2 = $d.Content
    $null = $headingRange2.Find.Execute("ОСНОВНЫЕ ФИНАНСОВЫЕ ПОКАЗАТЕЛИ КЛИЕНТА", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $headingPara2 = $headingRange2.Paragraphs(1)
    $headingPara2.Format.SpaceBefore = 0
}

# Edit 2: merge the two runs holding the "{issue.is_contract_corresponds_issuer_activity}"
# merge-field placeholder (previously split into "{issue.is_contract_corresponds_issuer_activity"
# and "}") back into a single run/run of text.
$fieldRange = $d.Content
$fieldRange.Find.Execute("{issue.is_contract_corresponds_issuer_activity}", $true, $false, $false, $false, $false, $true, 1, $false, "{issue.is_contract_corresponds_issuer_activity}", 2)
